$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1261
$ws.Range("F4").Value = 885
$ws.Range("F6").Value = 64
$ws.Range("F7").Value = 645
$ws.Range("F13").Value = 1309
$ws.Range("F16").Value = 523
$ws.Range("F17").Value = 741
$ws.Range("F18").Value = 29
$ws.Range("F19").Value = 278
$ws.Range("F22").Value = 10
$ws.Range("F24").Value = 4476
$ws.Range("F26").Value = 19
$ws.Range("F29").Value = 193
$ws.Range("F37").Value = 364
$ws.Range("F38").Value = 946
$ws.Range("F41").Value = 125

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1261
$ws.Range("F6").Value = 885
$ws.Range("F10").Value = 64
$ws.Range("F11").Value = 645
$ws.Range("F19").Value = 1309
$ws.Range("F22").Value = 523
$ws.Range("F24").Value = 741
$ws.Range("F25").Value = 29
$ws.Range("F26").Value = 278
$ws.Range("F28").Value = 10
$ws.Range("F29").Value = 4476
$ws.Range("F31").Value = 19
$ws.Range("F34").Value = 193
$ws.Range("F41").Value = 364
$ws.Range("F42").Value = 946
$ws.Range("F45").Value = 125
